$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 508-611: columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
$ws.Range("D508").Value = 44637
$ws.Range("J508").Value = 1650
$ws.Range("K508").Value = 1000
$ws.Range("L508").Value = 1050
$ws.Range("M508").Value = 1024
$ws.Range("O508").Value = "Provincia de Quillota"
$ws.Range("P508").Value = 1024

$ws.Range("D509").Value = 44637
$ws.Range("J509").Value = 880
$ws.Range("K509").Value = 800
$ws.Range("L509").Value = 800
$ws.Range("M509").Value = 800
$ws.Range("O509").Value = "Provincia de Quillota"
$ws.Range("P509").Value = 800

$ws.Range("D510").Value = 44208
$ws.Range("J510").Value = 1800
$ws.Range("K510").Value = 700
$ws.Range("L510").Value = 800
$ws.Range("M510").Value = 747
$ws.Range("O510").Value = "Provincia de Quillota"
$ws.Range("P510").Value = 747

$ws.Range("D511").Value = 44208
$ws.Range("J511").Value = 900
$ws.Range("K511").Value = 600
$ws.Range("L511").Value = 600
$ws.Range("M511").Value = 600
$ws.Range("O511").Value = "Provincia de Quillota"
$ws.Range("P511").Value = 600

$ws.Range("D512").Value = 44445
$ws.Range("J512").Value = 3000
$ws.Range("K512").Value = 550
$ws.Range("L512").Value = 600
$ws.Range("M512").Value = 580
$ws.Range("O512").Value = "Provincia de Quillota"
$ws.Range("P512").Value = 580

$ws.Range("D513").Value = 44445
$ws.Range("J513").Value = 1200
$ws.Range("K513").Value = 500
$ws.Range("L513").Value = 500
$ws.Range("M513").Value = 500
$ws.Range("O513").Value = "Provincia de Quillota"
$ws.Range("P513").Value = 500

$ws.Range("D514").Value = 44524
$ws.Range("J514").Value = 2550
$ws.Range("K514").Value = 500
$ws.Range("L514").Value = 550
$ws.Range("M514").Value = 526
$ws.Range("O514").Value = "Provincia de Quillota"
$ws.Range("P514").Value = 526

$ws.Range("D515").Value = 44524
$ws.Range("J515").Value = 1300
$ws.Range("K515").Value = 400
$ws.Range("L515").Value = 400
$ws.Range("M515").Value = 400
$ws.Range("O515").Value = "Provincia de Quillota"
$ws.Range("P515").Value = 400

$ws.Range("D516").Value = 44355
$ws.Range("J516").Value = 1300
$ws.Range("K516").Value = 700
$ws.Range("L516").Value = 700
$ws.Range("M516").Value = 700
$ws.Range("O516").Value = "Provincia de Quillota"
$ws.Range("P516").Value = 700

$ws.Range("D517").Value = 44355
$ws.Range("J517").Value = 880
$ws.Range("K517").Value = 600
$ws.Range("L517").Value = 600
$ws.Range("M517").Value = 600
$ws.Range("O517").Value = "Provincia de Quillota"
$ws.Range("P517").Value = 600

$ws.Range("D518").Value = 44530
$ws.Range("J518").Value = 2600
$ws.Range("K518").Value = 550
$ws.Range("L518").Value = 600
$ws.Range("M518").Value = 575
$ws.Range("O518").Value = "Provincia de Quillota"
$ws.Range("P518").Value = 575

$ws.Range("D519").Value = 44530
$ws.Range("J519").Value = 1200
$ws.Range("K519").Value = 400
$ws.Range("L519").Value = 400
$ws.Range("M519").Value = 400
$ws.Range("O519").Value = "Provincia de Quillota"
$ws.Range("P519").Value = 400

$ws.Range("D520").Value = 44483
$ws.Range("J520").Value = 4500
$ws.Range("K520").Value = 500
$ws.Range("L520").Value = 600
$ws.Range("M520").Value = 551
$ws.Range("O520").Value = "Provincia de Quillota"
$ws.Range("P520").Value = 551

$ws.Range("D521").Value = 44483
$ws.Range("J521").Value = 2200
$ws.Range("K521").Value = 400
$ws.Range("L521").Value = 400
$ws.Range("M521").Value = 400
$ws.Range("O521").Value = "Provincia de Quillota"
$ws.Range("P521").Value = 400

$ws.Range("D522").Value = 44294
$ws.Range("J522").Value = 900
$ws.Range("K522").Value = 800
$ws.Range("L522").Value = 800
$ws.Range("M522").Value = 800
$ws.Range("O522").Value = "Provincia de Quillota"
$ws.Range("P522").Value = 800

$ws.Range("D523").Value = 44294
$ws.Range("J523").Value = 950
$ws.Range("K523").Value = 700
$ws.Range("L523").Value = 700
$ws.Range("M523").Value = 700
$ws.Range("O523").Value = "Provincia de Quillota"
$ws.Range("P523").Value = 700

$ws.Range("D524").Value = 44617
$ws.Range("J524").Value = 600
$ws.Range("K524").Value = 1000
$ws.Range("L524").Value = 1000
$ws.Range("M524").Value = 1000
$ws.Range("O524").Value = "Provincia de Quillota"
$ws.Range("P524").Value = 1000

$ws.Range("D525").Value = 44617
$ws.Range("J525").Value = 550
$ws.Range("K525").Value = 800
$ws.Range("L525").Value = 800
$ws.Range("M525").Value = 800
$ws.Range("O525").Value = "Provincia de Quillota"
$ws.Range("P525").Value = 800

$ws.Range("D526").Value = 44557
$ws.Range("J526").Value = 2500
$ws.Range("K526").Value = 600
$ws.Range("L526").Value = 650
$ws.Range("M526").Value = 626
$ws.Range("O526").Value = "Provincia de Quillota"
$ws.Range("P526").Value = 626

$ws.Range("D527").Value = 44557
$ws.Range("J527").Value = 2200
$ws.Range("K527").Value = 450
$ws.Range("L527").Value = 500
$ws.Range("M527").Value = 477
$ws.Range("O527").Value = "Provincia de Quillota"
$ws.Range("P527").Value = 477

$ws.Range("D528").Value = 44489
$ws.Range("J528").Value = 1500
$ws.Range("K528").Value = 600
$ws.Range("L528").Value = 600
$ws.Range("M528").Value = 600
$ws.Range("O528").Value = "Provincia de Quillota"
$ws.Range("P528").Value = 600

$ws.Range("D529").Value = 44489
$ws.Range("J529").Value = 1600
$ws.Range("K529").Value = 500
$ws.Range("L529").Value = 500
$ws.Range("M529").Value = 500
$ws.Range("O529").Value = "Provincia de Quillota"
$ws.Range("P529").Value = 500

$ws.Range("D530").Value = 44264
$ws.Range("J530").Value = 950
$ws.Range("K530").Value = 800
$ws.Range("L530").Value = 800
$ws.Range("M530").Value = 800
$ws.Range("O530").Value = "Provincia de Quillota"
$ws.Range("P530").Value = 800

$ws.Range("D531").Value = 44264
$ws.Range("J531").Value = 900
$ws.Range("K531").Value = 700
$ws.Range("L531").Value = 700
$ws.Range("M531").Value = 700
$ws.Range("O531").Value = "Provincia de Quillota"
$ws.Range("P531").Value = 700

$ws.Range("D532").Value = 44396
$ws.Range("J532").Value = 1800
$ws.Range("K532").Value = 600
$ws.Range("L532").Value = 600
$ws.Range("M532").Value = 600
$ws.Range("O532").Value = "Provincia de Quillota"
$ws.Range("P532").Value = 600

$ws.Range("D533").Value = 44396
$ws.Range("J533").Value = 1600
$ws.Range("K533").Value = 500
$ws.Range("L533").Value = 500
$ws.Range("M533").Value = 500
$ws.Range("O533").Value = "Provincia de Quillota"
$ws.Range("P533").Value = 500

$ws.Range("D534").Value = 44232
$ws.Range("J534").Value = 900
$ws.Range("K534").Value = 900
$ws.Range("L534").Value = 900
$ws.Range("M534").Value = 900
$ws.Range("O534").Value = "Provincia de Quillota"
$ws.Range("P534").Value = 900

$ws.Range("D535").Value = 44232
$ws.Range("J535").Value = 880
$ws.Range("K535").Value = 700
$ws.Range("L535").Value = 700
$ws.Range("M535").Value = 700
$ws.Range("O535").Value = "Provincia de Quillota"
$ws.Range("P535").Value = 700

$ws.Range("D536").Value = 44279
$ws.Range("J536").Value = 950
$ws.Range("K536").Value = 900
$ws.Range("L536").Value = 900
$ws.Range("M536").Value = 900
$ws.Range("O536").Value = "Provincia de Santiago"
$ws.Range("P536").Value = 900

$ws.Range("D537").Value = 44279
$ws.Range("J537").Value = 880
$ws.Range("K537").Value = 700
$ws.Range("L537").Value = 700
$ws.Range("M537").Value = 700
$ws.Range("O537").Value = "Provincia de Santiago"
$ws.Range("P537").Value = 700

$ws.Range("D538").Value = 44330
$ws.Range("J538").Value = 1300
$ws.Range("K538").Value = 700
$ws.Range("L538").Value = 700
$ws.Range("M538").Value = 700
$ws.Range("O538").Value = "Provincia de Quillota"
$ws.Range("P538").Value = 700

$ws.Range("D539").Value = 44330
$ws.Range("J539").Value = 1280
$ws.Range("K539").Value = 600
$ws.Range("L539").Value = 600
$ws.Range("M539").Value = 600
$ws.Range("O539").Value = "Provincia de Quillota"
$ws.Range("P539").Value = 600

$ws.Range("D540").Value = 44504
$ws.Range("J540").Value = 2800
$ws.Range("K540").Value = 550
$ws.Range("L540").Value = 600
$ws.Range("M540").Value = 577
$ws.Range("O540").Value = "Provincia de Quillota"
$ws.Range("P540").Value = 577

$ws.Range("D541").Value = 44504
$ws.Range("J541").Value = 1500
$ws.Range("K541").Value = 400
$ws.Range("L541").Value = 400
$ws.Range("M541").Value = 400
$ws.Range("O541").Value = "Provincia de Quillota"
$ws.Range("P541").Value = 400

$ws.Range("D542").Value = 44572
$ws.Range("J542").Value = 2400
$ws.Range("K542").Value = 650
$ws.Range("L542").Value = 700
$ws.Range("M542").Value = 675
$ws.Range("O542").Value = "Provincia de Quillota"
$ws.Range("P542").Value = 675

$ws.Range("D543").Value = 44572
$ws.Range("J543").Value = 1600
$ws.Range("K543").Value = 500
$ws.Range("L543").Value = 500
$ws.Range("M543").Value = 500
$ws.Range("O543").Value = "Provincia de Quillota"
$ws.Range("P543").Value = 500

$ws.Range("D544").Value = 44257
$ws.Range("J544").Value = 1830
$ws.Range("K544").Value = 850
$ws.Range("L544").Value = 900
$ws.Range("M544").Value = 874
$ws.Range("O544").Value = "Provincia de Quillota"
$ws.Range("P544").Value = 874

$ws.Range("D545").Value = 44257
$ws.Range("J545").Value = 980
$ws.Range("K545").Value = 700
$ws.Range("L545").Value = 700
$ws.Range("M545").Value = 700
$ws.Range("O545").Value = "Provincia de Quillota"
$ws.Range("P545").Value = 700

$ws.Range("D546").Value = 44301
$ws.Range("J546").Value = 900
$ws.Range("K546").Value = 700
$ws.Range("L546").Value = 700
$ws.Range("M546").Value = 700
$ws.Range("O546").Value = "Provincia de Quillota"
$ws.Range("P546").Value = 700

$ws.Range("D547").Value = 44301
$ws.Range("J547").Value = 850
$ws.Range("K547").Value = 600
$ws.Range("L547").Value = 600
$ws.Range("M547").Value = 600
$ws.Range("O547").Value = "Provincia de Quillota"
$ws.Range("P547").Value = 600

$ws.Range("D548").Value = 44370
$ws.Range("J548").Value = 1900
$ws.Range("K548").Value = 600
$ws.Range("L548").Value = 600
$ws.Range("M548").Value = 600
$ws.Range("O548").Value = "Provincia de Quillota"
$ws.Range("P548").Value = 600

$ws.Range("D549").Value = 44370
$ws.Range("J549").Value = 1800
$ws.Range("K549").Value = 500
$ws.Range("L549").Value = 500
$ws.Range("M549").Value = 500
$ws.Range("O549").Value = "Provincia de Quillota"
$ws.Range("P549").Value = 500

$ws.Range("D550").Value = 44487
$ws.Range("J550").Value = 2150
$ws.Range("K550").Value = 550
$ws.Range("L550").Value = 600
$ws.Range("M550").Value = 578
$ws.Range("O550").Value = "Provincia de Quillota"
$ws.Range("P550").Value = 578

$ws.Range("D551").Value = 44487
$ws.Range("J551").Value = 1200
$ws.Range("K551").Value = 500
$ws.Range("L551").Value = 500
$ws.Range("M551").Value = 500
$ws.Range("O551").Value = "Provincia de Quillota"
$ws.Range("P551").Value = 500

$ws.Range("D552").Value = 44174
$ws.Range("J552").Value = 2700
$ws.Range("K552").Value = 600
$ws.Range("L552").Value = 650
$ws.Range("M552").Value = 628
$ws.Range("O552").Value = "Provincia de Quillota"
$ws.Range("P552").Value = 628

$ws.Range("D553").Value = 44174
$ws.Range("J553").Value = 1300
$ws.Range("K553").Value = 500
$ws.Range("L553").Value = 500
$ws.Range("M553").Value = 500
$ws.Range("O553").Value = "Provincia de Quillota"
$ws.Range("P553").Value = 500

$ws.Range("D554").Value = 44200
$ws.Range("J554").Value = 950
$ws.Range("K554").Value = 800
$ws.Range("L554").Value = 800
$ws.Range("M554").Value = 800
$ws.Range("O554").Value = "Provincia de Quillota"
$ws.Range("P554").Value = 800

$ws.Range("D555").Value = 44200
$ws.Range("J555").Value = 800
$ws.Range("K555").Value = 650
$ws.Range("L555").Value = 650
$ws.Range("M555").Value = 650
$ws.Range("O555").Value = "Provincia de Quillota"
$ws.Range("P555").Value = 650

$ws.Range("D556").Value = 44385
$ws.Range("J556").Value = 1800
$ws.Range("K556").Value = 600
$ws.Range("L556").Value = 600
$ws.Range("M556").Value = 600
$ws.Range("O556").Value = "Provincia de Quillota"
$ws.Range("P556").Value = 600

$ws.Range("D557").Value = 44385
$ws.Range("J557").Value = 1200
$ws.Range("K557").Value = 500
$ws.Range("L557").Value = 500
$ws.Range("M557").Value = 500
$ws.Range("O557").Value = "Provincia de Quillota"
$ws.Range("P557").Value = 500

$ws.Range("D558").Value = 44236
$ws.Range("J558").Value = 1950
$ws.Range("K558").Value = 950
$ws.Range("L558").Value = 1000
$ws.Range("M558").Value = 978
$ws.Range("O558").Value = "Provincia de Quillota"
$ws.Range("P558").Value = 978

$ws.Range("D559").Value = 44236
$ws.Range("J559").Value = 1200
$ws.Range("K559").Value = 800
$ws.Range("L559").Value = 800
$ws.Range("M559").Value = 800
$ws.Range("O559").Value = "Provincia de Quillota"
$ws.Range("P559").Value = 800

$ws.Range("D560").Value = 44221
$ws.Range("J560").Value = 1000
$ws.Range("K560").Value = 700
$ws.Range("L560").Value = 700
$ws.Range("M560").Value = 700
$ws.Range("O560").Value = "Provincia de Quillota"
$ws.Range("P560").Value = 700

$ws.Range("D561").Value = 44221
$ws.Range("J561").Value = 890
$ws.Range("K561").Value = 600
$ws.Range("L561").Value = 600
$ws.Range("M561").Value = 600
$ws.Range("O561").Value = "Provincia de Quillota"
$ws.Range("P561").Value = 600

$ws.Range("D562").Value = 44413
$ws.Range("J562").Value = 3400
$ws.Range("K562").Value = 600
$ws.Range("L562").Value = 650
$ws.Range("M562").Value = 626
$ws.Range("O562").Value = "Provincia de Quillota"
$ws.Range("P562").Value = 626

$ws.Range("D563").Value = 44413
$ws.Range("J563").Value = 1400
$ws.Range("K563").Value = 500
$ws.Range("L563").Value = 500
$ws.Range("M563").Value = 500
$ws.Range("O563").Value = "Provincia de Quillota"
$ws.Range("P563").Value = 500

$ws.Range("D564").Value = 44272
$ws.Range("J564").Value = 2600
$ws.Range("K564").Value = 800
$ws.Range("L564").Value = 850
$ws.Range("M564").Value = 827
$ws.Range("O564").Value = "Provincia de Quillota"
$ws.Range("P564").Value = 827

$ws.Range("D565").Value = 44272
$ws.Range("J565").Value = 1500
$ws.Range("K565").Value = 700
$ws.Range("L565").Value = 700
$ws.Range("M565").Value = 700
$ws.Range("O565").Value = "Provincia de Quillota"
$ws.Range("P565").Value = 700

$ws.Range("D566").Value = 44229
$ws.Range("J566").Value = 900
$ws.Range("K566").Value = 800
$ws.Range("L566").Value = 800
$ws.Range("M566").Value = 800
$ws.Range("O566").Value = "Provincia de Quillota"
$ws.Range("P566").Value = 800

$ws.Range("D567").Value = 44229
$ws.Range("J567").Value = 850
$ws.Range("K567").Value = 650
$ws.Range("L567").Value = 650
$ws.Range("M567").Value = 650
$ws.Range("O567").Value = "Provincia de Quillota"
$ws.Range("P567").Value = 650

$ws.Range("D568").Value = 44214
$ws.Range("J568").Value = 1600
$ws.Range("K568").Value = 700
$ws.Range("L568").Value = 700
$ws.Range("M568").Value = 700
$ws.Range("O568").Value = "Provincia de Quillota"
$ws.Range("P568").Value = 700

$ws.Range("D569").Value = 44214
$ws.Range("J569").Value = 1200
$ws.Range("K569").Value = 600
$ws.Range("L569").Value = 600
$ws.Range("M569").Value = 600
$ws.Range("O569").Value = "Provincia de Quillota"
$ws.Range("P569").Value = 600

$ws.Range("D570").Value = 44299
$ws.Range("J570").Value = 1100
$ws.Range("K570").Value = 800
$ws.Range("L570").Value = 800
$ws.Range("M570").Value = 800
$ws.Range("O570").Value = "Provincia de Quillota"
$ws.Range("P570").Value = 800

$ws.Range("D571").Value = 44299
$ws.Range("J571").Value = 880
$ws.Range("K571").Value = 700
$ws.Range("L571").Value = 700
$ws.Range("M571").Value = 700
$ws.Range("O571").Value = "Provincia de Quillota"
$ws.Range("P571").Value = 700

$ws.Range("D572").Value = 44610
$ws.Range("J572").Value = 850
$ws.Range("K572").Value = 1000
$ws.Range("L572").Value = 1000
$ws.Range("M572").Value = 1000
$ws.Range("O572").Value = "Provincia de Quillota"
$ws.Range("P572").Value = 1000

$ws.Range("D573").Value = 44610
$ws.Range("J573").Value = 880
$ws.Range("K573").Value = 800
$ws.Range("L573").Value = 800
$ws.Range("M573").Value = 800
$ws.Range("O573").Value = "Provincia de Quillota"
$ws.Range("P573").Value = 800

$ws.Range("D574").Value = 44312
$ws.Range("J574").Value = 1300
$ws.Range("K574").Value = 800
$ws.Range("L574").Value = 800
$ws.Range("M574").Value = 800
$ws.Range("O574").Value = "Provincia de Quillota"
$ws.Range("P574").Value = 800

$ws.Range("D575").Value = 44312
$ws.Range("J575").Value = 950
$ws.Range("K575").Value = 700
$ws.Range("L575").Value = 700
$ws.Range("M575").Value = 700
$ws.Range("O575").Value = "Provincia de Quillota"
$ws.Range("P575").Value = 700

$ws.Range("D576").Value = 44399
$ws.Range("J576").Value = 2400
$ws.Range("K576").Value = 600
$ws.Range("L576").Value = 650
$ws.Range("M576").Value = 625
$ws.Range("O576").Value = "Provincia de Quillota"
$ws.Range("P576").Value = 625

$ws.Range("D577").Value = 44399
$ws.Range("J577").Value = 1300
$ws.Range("K577").Value = 500
$ws.Range("L577").Value = 500
$ws.Range("M577").Value = 500
$ws.Range("O577").Value = "Provincia de Quillota"
$ws.Range("P577").Value = 500

$ws.Range("D578").Value = 44615
$ws.Range("J578").Value = 600
$ws.Range("K578").Value = 1000
$ws.Range("L578").Value = 1000
$ws.Range("M578").Value = 1000
$ws.Range("O578").Value = "Provincia de Quillota"
$ws.Range("P578").Value = 1000

$ws.Range("D579").Value = 44615
$ws.Range("J579").Value = 680
$ws.Range("K579").Value = 800
$ws.Range("L579").Value = 800
$ws.Range("M579").Value = 800
$ws.Range("O579").Value = "Provincia de Quillota"
$ws.Range("P579").Value = 800

$ws.Range("D580").Value = 44522
$ws.Range("J580").Value = 2600
$ws.Range("K580").Value = 500
$ws.Range("L580").Value = 550
$ws.Range("M580").Value = 523
$ws.Range("O580").Value = "Provincia de Quillota"
$ws.Range("P580").Value = 523

$ws.Range("D581").Value = 44522
$ws.Range("J581").Value = 1400
$ws.Range("K581").Value = 400
$ws.Range("L581").Value = 400
$ws.Range("M581").Value = 400
$ws.Range("O581").Value = "Provincia de Quillota"
$ws.Range("P581").Value = 400

$ws.Range("D582").Value = 44543
$ws.Range("J582").Value = 1300
$ws.Range("K582").Value = 600
$ws.Range("L582").Value = 600
$ws.Range("M582").Value = 600
$ws.Range("O582").Value = "Provincia de Quillota"
$ws.Range("P582").Value = 600

$ws.Range("D583").Value = 44543
$ws.Range("J583").Value = 2700
$ws.Range("K583").Value = 400
$ws.Range("L583").Value = 500
$ws.Range("M583").Value = 452
$ws.Range("O583").Value = "Provincia de Quillota"
$ws.Range("P583").Value = 452

$ws.Range("D584").Value = 44167
$ws.Range("J584").Value = 1800
$ws.Range("K584").Value = 600
$ws.Range("L584").Value = 650
$ws.Range("M584").Value = 626
$ws.Range("O584").Value = "Provincia de Quillota"
$ws.Range("P584").Value = 626

$ws.Range("D585").Value = 44167
$ws.Range("J585").Value = 950
$ws.Range("K585").Value = 500
$ws.Range("L585").Value = 500
$ws.Range("M585").Value = 500
$ws.Range("O585").Value = "Provincia de Quillota"
$ws.Range("P585").Value = 500

$ws.Range("D586").Value = 44277
$ws.Range("J586").Value = 2100
$ws.Range("K586").Value = 800
$ws.Range("L586").Value = 900
$ws.Range("M586").Value = 852
$ws.Range("O586").Value = "Provincia de Quillota"
$ws.Range("P586").Value = 852

$ws.Range("D587").Value = 44277
$ws.Range("J587").Value = 1200
$ws.Range("K587").Value = 700
$ws.Range("L587").Value = 700
$ws.Range("M587").Value = 700
$ws.Range("O587").Value = "Provincia de Quillota"
$ws.Range("P587").Value = 700

$ws.Range("D588").Value = 44258
$ws.Range("J588").Value = 1100
$ws.Range("K588").Value = 800
$ws.Range("L588").Value = 800
$ws.Range("M588").Value = 800
$ws.Range("O588").Value = "Provincia de Quillota"
$ws.Range("P588").Value = 800

$ws.Range("D589").Value = 44258
$ws.Range("J589").Value = 950
$ws.Range("K589").Value = 700
$ws.Range("L589").Value = 700
$ws.Range("M589").Value = 700
$ws.Range("O589").Value = "Provincia de Quillota"
$ws.Range("P589").Value = 700

$ws.Range("D590").Value = 44390
$ws.Range("J590").Value = 1700
$ws.Range("K590").Value = 600
$ws.Range("L590").Value = 600
$ws.Range("M590").Value = 600
$ws.Range("O590").Value = "Provincia de Quillota"
$ws.Range("P590").Value = 600

$ws.Range("D591").Value = 44390
$ws.Range("J591").Value = 2280
$ws.Range("K591").Value = 450
$ws.Range("L591").Value = 500
$ws.Range("M591").Value = 471
$ws.Range("O591").Value = "Provincia de Quillota"
$ws.Range("P591").Value = 471

$ws.Range("D592").Value = 44349
$ws.Range("J592").Value = 1800
$ws.Range("K592").Value = 700
$ws.Range("L592").Value = 700
$ws.Range("M592").Value = 700
$ws.Range("O592").Value = "Provincia de Quillota"
$ws.Range("P592").Value = 700

$ws.Range("D593").Value = 44349
$ws.Range("J593").Value = 1600
$ws.Range("K593").Value = 600
$ws.Range("L593").Value = 600
$ws.Range("M593").Value = 600
$ws.Range("O593").Value = "Provincia de Quillota"
$ws.Range("P593").Value = 600

$ws.Range("D594").Value = 44285
$ws.Range("J594").Value = 1900
$ws.Range("K594").Value = 750
$ws.Range("L594").Value = 800
$ws.Range("M594").Value = 775
$ws.Range("O594").Value = "Provincia de Quillota"
$ws.Range("P594").Value = 775

$ws.Range("D595").Value = 44285
$ws.Range("J595").Value = 980
$ws.Range("K595").Value = 600
$ws.Range("L595").Value = 600
$ws.Range("M595").Value = 600
$ws.Range("O595").Value = "Provincia de Quillota"
$ws.Range("P595").Value = 600

$ws.Range("D596").Value = 44498
$ws.Range("J596").Value = 2500
$ws.Range("K596").Value = 600
$ws.Range("L596").Value = 650
$ws.Range("M596").Value = 624
$ws.Range("O596").Value = "Provincia de Melipilla"
$ws.Range("P596").Value = 624

$ws.Range("D597").Value = 44498
$ws.Range("J597").Value = 1400
$ws.Range("K597").Value = 500
$ws.Range("L597").Value = 500
$ws.Range("M597").Value = 500
$ws.Range("O597").Value = "Provincia de Melipilla"
$ws.Range("P597").Value = 500

$ws.Range("D598").Value = 44179
$ws.Range("J598").Value = 1720
$ws.Range("K598").Value = 650
$ws.Range("L598").Value = 700
$ws.Range("M598").Value = 675
$ws.Range("O598").Value = "Provincia de Quillota"
$ws.Range("P598").Value = 675

$ws.Range("D599").Value = 44179
$ws.Range("J599").Value = 1580
$ws.Range("K599").Value = 500
$ws.Range("L599").Value = 550
$ws.Range("M599").Value = 525
$ws.Range("O599").Value = "Provincia de Quillota"
$ws.Range("P599").Value = 525

$ws.Range("D600").Value = 44418
$ws.Range("J600").Value = 3300
$ws.Range("K600").Value = 600
$ws.Range("L600").Value = 650
$ws.Range("M600").Value = 623
$ws.Range("O600").Value = "Provincia de Quillota"
$ws.Range("P600").Value = 623

$ws.Range("D601").Value = 44418
$ws.Range("J601").Value = 1600
$ws.Range("K601").Value = 500
$ws.Range("L601").Value = 500
$ws.Range("M601").Value = 500
$ws.Range("O601").Value = "Provincia de Quillota"
$ws.Range("P601").Value = 500

$ws.Range("D602").Value = 44595
$ws.Range("J602").Value = 850
$ws.Range("K602").Value = 900
$ws.Range("L602").Value = 900
$ws.Range("M602").Value = 900
$ws.Range("O602").Value = "Provincia de Quillota"
$ws.Range("P602").Value = 900

$ws.Range("D603").Value = 44595
$ws.Range("J603").Value = 900
$ws.Range("K603").Value = 700
$ws.Range("L603").Value = 700
$ws.Range("M603").Value = 700
$ws.Range("O603").Value = "Provincia de Quillota"
$ws.Range("P603").Value = 700

$ws.Range("D604").Value = 44628
$ws.Range("J604").Value = 1300
$ws.Range("K604").Value = 1000
$ws.Range("L604").Value = 1100
$ws.Range("M604").Value = 1050
$ws.Range("O604").Value = "Provincia de Quillota"
$ws.Range("P604").Value = 1050

$ws.Range("D605").Value = 44628
$ws.Range("J605").Value = 600
$ws.Range("K605").Value = 800
$ws.Range("L605").Value = 800
$ws.Range("M605").Value = 800
$ws.Range("O605").Value = "Provincia de Quillota"
$ws.Range("P605").Value = 800

$ws.Range("D606").Value = 44335
$ws.Range("J606").Value = 1200
$ws.Range("K606").Value = 700
$ws.Range("L606").Value = 700
$ws.Range("M606").Value = 700
$ws.Range("O606").Value = "Provincia de Quillota"
$ws.Range("P606").Value = 700

$ws.Range("D607").Value = 44335
$ws.Range("J607").Value = 1100
$ws.Range("K607").Value = 600
$ws.Range("L607").Value = 600
$ws.Range("M607").Value = 600
$ws.Range("O607").Value = "Provincia de Quillota"
$ws.Range("P607").Value = 600

$ws.Range("D608").Value = 44552
$ws.Range("J608").Value = 1800
$ws.Range("K608").Value = 600
$ws.Range("L608").Value = 600
$ws.Range("M608").Value = 600
$ws.Range("O608").Value = "Provincia de Quillota"
$ws.Range("P608").Value = 600

$ws.Range("D609").Value = 44552
$ws.Range("J609").Value = 3100
$ws.Range("K609").Value = 450
$ws.Range("L609").Value = 500
$ws.Range("M609").Value = 481
$ws.Range("O609").Value = "Provincia de Quillota"
$ws.Range("P609").Value = 481

$ws.Range("D610").Value = 44544
$ws.Range("J610").Value = 2450
$ws.Range("K610").Value = 550
$ws.Range("L610").Value = 600
$ws.Range("M610").Value = 574
$ws.Range("O610").Value = "Provincia de Quillota"
$ws.Range("P610").Value = 574

$ws.Range("D611").Value = 44544
$ws.Range("J611").Value = 1300
$ws.Range("K611").Value = 450
$ws.Range("L611").Value = 450
$ws.Range("M611").Value = 450
$ws.Range("O611").Value = "Provincia de Quillota"
$ws.Range("P611").Value = 450

# Add new rows 612 and 613 (new data appended at the bottom of the table)
$ws.Range("A612").Value = 3
$ws.Range("B612").Value = "Femacal de La Calera"
$ws.Range("C612").Value = "Coquimbo"
$ws.Range("D612").Value = 44160
$ws.Range("E612").Value = 5
$ws.Range("F612").Value = 100112023
$ws.Range("G612").Value = "Brócoli"
$ws.Range("H612").Value = "Sin especificar"
$ws.Range("I612").Value = "Primera"
$ws.Range("J612").Value = 980
$ws.Range("K612").Value = 600
$ws.Range("L612").Value = 600
$ws.Range("M612").Value = 600
$ws.Range("N612").Value = "$/unidad"
$ws.Range("O612").Value = "Provincia de Quillota"
$ws.Range("P612").Value = 600
$ws.Range("Q612").Value = 1
$ws.Range("R612").Value = "Hortaliza"

$ws.Range("A613").Value = 3
$ws.Range("B613").Value = "Femacal de La Calera"
$ws.Range("C613").Value = "Coquimbo"
$ws.Range("D613").Value = 44160
$ws.Range("E613").Value = 5
$ws.Range("F613").Value = 100112023
$ws.Range("G613").Value = "Brócoli"
$ws.Range("H613").Value = "Sin especificar"
$ws.Range("I613").Value = "Segunda"
$ws.Range("J613").Value = 900
$ws.Range("K613").Value = 500
$ws.Range("L613").Value = 500
$ws.Range("M613").Value = 500
$ws.Range("N613").Value = "$/unidad"
$ws.Range("O613").Value = "Provincia de Quillota"
$ws.Range("P613").Value = 500
$ws.Range("Q613").Value = 1
$ws.Range("R613").Value = "Hortaliza"

# Update sheet dimension to reflect the new used range
$ws.Range("D612").NumberFormat = $ws.Range("D611").NumberFormat
$ws.Range("D613").NumberFormat = $ws.Range("D611").NumberFormat